$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48, shifting existing rows 48:55 down to 49:56
$ws.Rows.Item(48).Insert(4)

# Populate the new row 48 with the new weekly record.
# Columns A,B,C,E,F,G,H,N,O,Q,R stay identical to the neighbouring records.
$ws.Cells.Item(48, 1).Value = 1
$ws.Cells.Item(48, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value = 44943
$ws.Cells.Item(48, 5).Value = 15
$ws.Cells.Item(48, 6).Value = 100112028
$ws.Cells.Item(48, 7).Value = "Sandia"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Especial"
$ws.Cells.Item(48, 10).Value = 800
$ws.Cells.Item(48, 11).Value = 750
$ws.Cells.Item(48, 12).Value = 780
$ws.Cells.Item(48, 13).Value = 769
$ws.Cells.Item(48, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(48, 15).Value = "Perú"
$ws.Cells.Item(48, 16).Value = 769
$ws.Cells.Item(48, 17).Value = 1
$ws.Cells.Item(48, 18).Value = "Hortaliza"
